$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '245.40'

Set-TextValue "D3" '23.85'

Set-TextValue "D4" '5.384'

Set-TextValue "D5" '0.05896'

$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D6" '3.383'
$ws.Range("E6").Value = '5GateTokenGT'

$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D7" '6.493'
$ws.Range("E7").Value = '6KuCoinTokenKCS'

Set-TextValue "D8" '0.8100'

Set-TextValue "D9" '0.9243'

Set-TextValue "D10" '0.1423'

Set-TextValue "D11" '0.07418'

Set-TextValue "D12" '0.03129'

Set-TextValue "D13" '0.03044'

Set-TextValue "D14" '0.09362'

Set-TextValue "D15" '3.849'

Set-TextValue "D16" '0.001568'

Set-TextValue "D17" '0.04703'

$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue "D18" '0.0005990'
$ws.Range("E18").Value = '17OneONEWorstin24h'

$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D19" '0.005924'
$ws.Range("E19").Value = '18TigerCashTCH'

$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue "D20" '0.001241'
$ws.Range("E20").Value = '19BitKanKAN'

$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue "D21" '0.004720'
$ws.Range("E21").Value = '20HotbitTokenHTB'

$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue "D22" '0.00008803'
$ws.Range("E22").Value = '21NitroExNTXBestin24h'

$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D23" '3.562'
$ws.Range("E23").Value = '22LEOLEO'

$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D24" '2.158'
$ws.Range("E24").Value = '23BTSETokenBTSE'

Set-TextValue "D25" '0.3229'

Set-TextValue "D27" '0.0002653'

Set-TextValue "D40" '0.03890'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D41" '0.1068'
$ws.Range("E41").Value = '40BKEXTokenBKK'

Set-TextValue "D42" '0.002761'

$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D43" '0.003105'
$ws.Range("E43").Value = '42KickTokenKICK'

Set-TextValue "D44" '0.008525'

Set-TextValue "D45" '0.00005251'

Set-TextValue "D47" '0.6710'

Set-TextValue "D48" '0.001942'
$ws.Range("E48").Value = '47BOLOBOLO'

Set-TextValue "D49" '0.00002100'

Set-TextValue "D50" '0.0002000'
